# Generate Report for Handoff
#
# The "b.md" row (row 3) in each sheet moves from
# "Handed back: in sync with en-US" to "Ready for handoff", a new handoff
# xliff file is generated for each locale, and an error detail message is
# recorded explaining that the handback file is stale.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-21 08:44:08"

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# A literal "False" gets auto-typed as a Boolean by plain .Value assignment;
# a leading apostrophe forces text entry (and is stripped from the stored
# text), then resetting the style to Normal drops the quote-prefix marker
# so the cell keeps its original (default) style, same as the rest of the
# "True"/"False" text cells in this sheet.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-21 08:44:00"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85d8ab2d5b7358e9f79f15655dd4cd528c9fbc86/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99db309d3db9dfc004407604a0d7fe881d18f76c/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.1

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-21 08:44:08"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85d8ab2d5b7358e9f79f15655dd4cd528c9fbc86/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99db309d3db9dfc004407604a0d7fe881d18f76c/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.1
